$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.0000004157731461785933
$ws.Range("E3").Value = 0.0000004157731461785933

$ws.Range("D4").Value = 0.957011486874055
$ws.Range("E4").Value = 0.957011486874055

$ws.Range("D5").Value = 0.00002033956699947414
$ws.Range("E5").Value = 0.00002033956699947414

$ws.Range("D6").Value = 0.0000000000000000000000007064351745193415
$ws.Range("E6").Value = 0.0000000000000000000000007064351745193415

$ws.Range("D7").Value = 0.005965608981428
$ws.Range("E7").Value = 0.994034391018572

$ws.Range("D8").Value = 0.00002269328663683384
$ws.Range("E8").Value = 0.9999773067133632

$ws.Range("D9").Value = 0.9999999999965707
$ws.Range("E9").Value = 0.000000000003429256878462184

$ws.Range("D10").Value = 0.000005019302784133008
$ws.Range("E10").Value = 0.9999949806972158

$ws.Range("D11").Value = 0.9999673088640366
$ws.Range("E11").Value = 0.0000326911359633808
$ws.Range("F11").Value = 9.984767913818359
